$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DPLKINV161-001")
$ws2 = $wb.Worksheets.Item("DPLKINV161-002")

# ---------------------------------------------------------------------------
# Sheet DPLKINV161-001 ("Disetujui" / approve branch)
#   - the switching id used to be the generic "Hasil Generate" placeholder;
#     it is now filled in with the real generated id.
# ---------------------------------------------------------------------------
$ws1.Range("N2").Value = "SWRKS202200008"
$ws1.Range("F2").Value = "Username : 33372;`r`nPassword : bni1234;`r`nRole : 18 - Pimpinan Kelompok Investasi;`r`nDeviden ID : SWRKS202200008;`r`nStatus Verifikasi : 1 : Setuju;`r`nKeterangan Verifikasi : INV.RKS.TRX.004 PEMBATALAN DISETUJUI"

# ---------------------------------------------------------------------------
# Sheet DPLKINV161-002 ("Kembalikan ke Data Entry" / return branch)
#   - N2 was blank, now also carries the generated switching id.
# ---------------------------------------------------------------------------
$ws2.Range("N2").Value = "SWRKS202200007"
$ws2.Range("F2").Value = "Username : 33372;`r`nPassword : bni1234;`r`nRole : 18 - Pimpinan Kelompok Investasi;`r`nSwitching ID : SWRKS202200007;`r`nStatus Verifikasi : 0 : Kembalikan ke Data Entry;`r`nKeterangan Verifikasi : INV.RKS.TRX.004 Pembatalan dikembalikan ke Data Entry"

# ---------------------------------------------------------------------------
# Column N got a little wider on both sheets to fit "SWRKS202200008".
# ---------------------------------------------------------------------------
$ws1.Columns.Item(14).ColumnWidth = 14.67
$ws2.Columns.Item(14).ColumnWidth = 14.67

# Row 2 on sheet 1 grew taller (manual resize while reviewing the new text).
$ws1.Rows.Item(2).RowHeight = 127.5

# ---------------------------------------------------------------------------
# View state: the author ended up on DPLKINV161-001 with I2 selected,
# leaving DPLKINV161-002's own selection untouched at O2.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("I2").Select()
